# Commit: "Change scope to range for clarity"
#
# The "Parameters" column (column E) of the "Tabelle1" sheet contains a
# number of entries describing numeric bounds, written like "scope=0-300".
# This renames the "scope=" prefix to "range=" for every such entry.
# (Shared strings are de-duplicated under the hood, so several rows that
# happen to share an identical "scope=..." value are updated together by
# writing each referencing cell.)
#
# It also restores the saved view/selection back to the top of the sheet
# (topLeftCell A1, active cell/selection E19) instead of the scrolled-down
# state (topLeftCell A64, active cell/selection E83) that was saved before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E3").Value  = "range=0-300"
$ws.Range("E4").Value  = "range=0-300"
$ws.Range("E5").Value  = "range=01067-99998"
$ws.Range("E6").Value  = "range=100000000-999999999"
$ws.Range("E7").Value  = "range=1000000000-9999999999"
$ws.Range("E10").Value = "range=0.0-45.0"
$ws.Range("E12").Value = "range=0-1"
$ws.Range("E13").Value = "range=1-5"
$ws.Range("E15").Value = "range=1000000000-9999999999"
$ws.Range("E18").Value = "range=0-99"
$ws.Range("E19").Value = "range=0-6"
$ws.Range("E20").Value = "range=1-4"
$ws.Range("E35").Value = "range=1-6"
$ws.Range("E36").Value = "range=1-5"
$ws.Range("E37").Value = "range=0-99"
$ws.Range("E38").Value = "range=1-6"
$ws.Range("E39").Value = "range=0-100"
$ws.Range("E40").Value = "range=0-10"
$ws.Range("E41").Value = "range=3-15"
$ws.Range("E73").Value = "range=10-120"
$ws.Range("E74").Value = "range=10-120"
$ws.Range("E75").Value = "range=10-120"
$ws.Range("E76").Value = "range=10-120"
$ws.Range("E77").Value = "range=10-120"
$ws.Range("E78").Value = "range=10-120"

# Reset the saved scroll position / selection: topLeftCell back to A1, and
# the active cell/selection back to E19.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E19").Select()
